$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

$ws.Range("A2").Value = "A. From tree rings back to:1400 AD"
$ws.Range("F2").Value = "Cap depletions to 10-year lookback:period of flow."
$ws.Range("A3").Value = "B. Collaborator choices in:immersive modeling:sessions"
$ws.Range("F3").Value = "Divide inflow; Users consume and:conserve within their:account balance."
$ws.Range("A4").Value = "C. 85%, 65%, and 50% of:2000 to 2018:average flow"
$ws.Range("F4").Value = "Release 95% of regulated:inflow."
$ws.Range("A5").Value = "D. Reclamation's:Post 2026:web tool"
$ws.Range("F5").Value = "Release to prevent drawdown:to 3,490 feet."
$ws.Range("A6").Value = "E. Low Lake Powell releases:+ gains through:Grand Canyon"
$ws.Range("A7").Value = "F. Lowest consecutive flows:in Reclamation's ensembles:and traces"
$ws.Range("F6").Value = "Rule curve; Consumption equals or:less than inflow minus:evaporation."

$ws.Range("F7").Select() | Out-Null
